$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "104 capacitor code"
$ws.Range("E5").Value = "103 capacitor code"
$ws.Range("D4").Value = "Electrolytic Capacitor"

$ws.Range("E4").Select()

$wb.Save()
